# #67 Ajouter l'année 2020 dans la fiche de stats
# Corrige les libellés d'années des en-têtes de colonnes (ligne 3) sur les
# deux feuilles "Par saison" pour qu'ils soient cohérents avec les données
# (colonnes D:F = 2020, G:I = 2019, J:L = 2018, M:O = 2017).

$wb = $excel.ActiveWorkbook

$sheetNames = @("Par saison (fin de saison)", "Par saison (date de génération)")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("D3").Value = "2020-2019"
    $ws.Range("J3").Value = "2018-2017"
    $ws.Range("M3").Value = "2017-2016"
}
